# Add new word / "seen" + "imgPath" tracking columns to the fruit vocab sheet,
# fix a typo in the last Vietnamese translation, and mark the first row ("apple")
# as already "seen".

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- New header cells: C2 = "seen", D2 = "imgPath" (bold, like the other headers) ---
$ws.Range("C2").Value = "seen"
$ws.Range("C2").Font.Bold = $true

$ws.Range("D2").Value = "imgPath"
$ws.Range("D2").Font.Bold = $true

# --- New "seen" column: TRUE for the first word (apple), FALSE for the rest ---
$seenRange = $ws.Range("C3:C19")
$seenRange.NumberFormat = """TRUE"";""TRUE"";""FALSE"""

$ws.Range("C3").Formula = "=TRUE()"

for ($row = 4; $row -le 19; $row++) {
    $ws.Cells.Item($row, 3).Formula = "=FALSE()"
}

# --- Fix typo: last row's Vietnamese placeholder "clgt??" -> "hugh??" ---
$ws.Range("B19").Value = "hugh??"

# --- Cosmetic: column widths shrank slightly and the selection moved to C3 ---
$ws.Columns.Item(1).ColumnWidth = 14.0
$ws.Columns.Item(2).ColumnWidth = 24.67
$ws.Columns.Item(3).ColumnWidth = 11.62

$ws.Range("C3").Select() | Out-Null
